$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "data_CCM-code_FLASH"
$wb.Worksheets.Item(2).Name = "data_CCM-code_CCM"
$wb.Worksheets.Item(3).Name = "data_RAM-code_FLASH"
$wb.Worksheets.Item(4).Name = "data_RAM-code_CCM"
